# "Turn deref into a factory"
#
# Adds a new "no_deref" row of sample data to the TestSheet, right between
# the existing "wealth" row (row 3) and the trailing row 5, and moves the
# active selection down to B6 to reflect the sheet's new extent.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 4: label "no_deref" in column A, value 50 in column B.
$ws.Range("A4").Value = "no_deref"
$ws.Range("B4").Value = 50

# Update the selected cell to B6, matching the sheet's saved view state.
$ws.Range("B6").Select()
